$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (Förändrad) holds a date serial number that needs to be
# advanced from 45683 (2025-01-26) to 45684 (2025-01-27) for every
# data row (rows 2 through 36).
$ws.Range("C2:C36").Value = 45684
